$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 178
$ws.Range("B178").Value = 6937270
$ws.Range("E178").Value = 'Olympiakos'
$ws.Range("F178").Value = 'Volos NFC'
$ws.Range("G178").Value = 3
$ws.Range("H178").Value = 0
$ws.Range("I178").Value = 2
$ws.Range("J178").Value = 0
$ws.Range("K178").Value = 'H'
$ws.Range("L178").Value = 1.125
$ws.Range("M178").Value = 9
$ws.Range("N178").Value = 19
$ws.Range("O178").Value = 1.111
$ws.Range("P178").Value = 9
$ws.Range("Q178").Value = 21
$ws.Range("R178").Value = -2.25
$ws.Range("S178").Value = 1.875
$ws.Range("T178").Value = 1.975
$ws.Range("U178").Value = 3.25
$ws.Range("V178").Value = 2
$ws.Range("W178").Value = 1.85
$ws.Range("X178").Value = 0.111
$ws.Range("Y178").Value = -1
$ws.Range("Z178").Value = -1
$ws.Range("AA178").Value = 0.875
$ws.Range("AB178").Value = -1
$ws.Range("AC178").Value = -0.5
$ws.Range("AD178").Value = 0.425

# Row 179
$ws.Range("B179").Value = 6937271
$ws.Range("E179").Value = 'Giannina'
$ws.Range("F179").Value = 'Atromitos Athinon'
$ws.Range("G179").Value = 1
$ws.Range("H179").Value = 1
$ws.Range("I179").Value = 1
$ws.Range("J179").Value = 0
$ws.Range("K179").Value = 'D'
$ws.Range("L179").Value = 2.45
$ws.Range("M179").Value = 3.1
$ws.Range("N179").Value = 3.1
$ws.Range("O179").Value = 2
$ws.Range("P179").Value = 3.3
$ws.Range("Q179").Value = 4
$ws.Range("R179").Value = -0.5
$ws.Range("S179").Value = 2.025
$ws.Range("T179").Value = 1.825
$ws.Range("U179").Value = 2.25
$ws.Range("V179").Value = 1.85
$ws.Range("W179").Value = 2
$ws.Range("X179").Value = -1
$ws.Range("Y179").Value = 2.3
$ws.Range("Z179").Value = -1
$ws.Range("AA179").Value = -1
$ws.Range("AB179").Value = 0.825
$ws.Range("AC179").Value = -0.5
$ws.Range("AD179").Value = 0.5

# Row 180
$ws.Range("B180").Value = 6937272
$ws.Range("E180").Value = 'Lamia'
$ws.Range("F180").Value = 'PAOK Salonika'
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 2
$ws.Range("I180").Value = 0
$ws.Range("J180").Value = 0
$ws.Range("K180").Value = 'A'
$ws.Range("L180").Value = 7.5
$ws.Range("M180").Value = 4.5
$ws.Range("N180").Value = 1.444
$ws.Range("O180").Value = 9.5
$ws.Range("P180").Value = 5
$ws.Range("Q180").Value = 1.333
$ws.Range("R180").Value = 1.5
$ws.Range("S180").Value = 1.925
$ws.Range("T180").Value = 1.925
$ws.Range("U180").Value = 3
$ws.Range("V180").Value = 1.95
$ws.Range("W180").Value = 1.9
$ws.Range("X180").Value = -1
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = 0.333
$ws.Range("AA180").Value = -1
$ws.Range("AB180").Value = 0.925
$ws.Range("AC180").Value = -1
$ws.Range("AD180").Value = 0.8999999999999999

# Row 181
$ws.Range("B181").Value = 6937269
$ws.Range("E181").Value = 'Aris Salonika'
$ws.Range("F181").Value = 'AEK Athens'
$ws.Range("G181").Value = 3
$ws.Range("H181").Value = 3
$ws.Range("I181").Value = 1
$ws.Range("J181").Value = 1
$ws.Range("K181").Value = 'D'
$ws.Range("L181").Value = 4.75
$ws.Range("M181").Value = 3.75
$ws.Range("N181").Value = 1.75
$ws.Range("O181").Value = 6.5
$ws.Range("P181").Value = 4.2
$ws.Range("Q181").Value = 1.5
$ws.Range("R181").Value = 1
$ws.Range("S181").Value = 2.05
$ws.Range("T181").Value = 1.8
$ws.Range("U181").Value = 2.5
$ws.Range("V181").Value = 1.975
$ws.Range("W181").Value = 1.875
$ws.Range("X181").Value = -1
$ws.Range("Y181").Value = 3.2
$ws.Range("Z181").Value = -1
$ws.Range("AA181").Value = 1.05
$ws.Range("AB181").Value = -1
$ws.Range("AC181").Value = 0.9750000000000001
$ws.Range("AD181").Value = -1

# Row 194
$ws.Range("B194").Value = 7920470
$ws.Range("E194").Value = 'AEK Athens'
$ws.Range("F194").Value = 'Olympiakos'
$ws.Range("G194").Value = 1
$ws.Range("H194").Value = 0
$ws.Range("I194").Value = 0
$ws.Range("J194").Value = 0
$ws.Range("K194").Value = 'H'
$ws.Range("L194").Value = 1.909
$ws.Range("M194").Value = 3.4
$ws.Range("N194").Value = 4.2
$ws.Range("O194").Value = 2.2
$ws.Range("P194").Value = 3.2
$ws.Range("Q194").Value = 3.5
$ws.Range("R194").Value = -0.25
$ws.Range("S194").Value = 1.85
$ws.Range("T194").Value = 2
$ws.Range("U194").Value = 2.5
$ws.Range("V194").Value = 2.025
$ws.Range("W194").Value = 1.825
$ws.Range("X194").Value = 1.2
$ws.Range("Y194").Value = -1
$ws.Range("Z194").Value = -1
$ws.Range("AA194").Value = 0.8500000000000001
$ws.Range("AB194").Value = -1
$ws.Range("AC194").Value = -1
$ws.Range("AD194").Value = 0.825

# Row 195
$ws.Range("B195").Value = 7920471
$ws.Range("E195").Value = 'Aris Salonika'
$ws.Range("F195").Value = 'Lamia'
$ws.Range("G195").Value = 3
$ws.Range("H195").Value = 1
$ws.Range("I195").Value = 0
$ws.Range("J195").Value = 0
$ws.Range("K195").Value = 'H'
$ws.Range("L195").Value = 1.571
$ws.Range("M195").Value = 4
$ws.Range("N195").Value = 6
$ws.Range("O195").Value = 1.444
$ws.Range("P195").Value = 4.5
$ws.Range("Q195").Value = 8.5
$ws.Range("R195").Value = -1.25
$ws.Range("S195").Value = 1.925
$ws.Range("T195").Value = 1.925
$ws.Range("U195").Value = 2.75
$ws.Range("V195").Value = 2.025
$ws.Range("W195").Value = 1.825
$ws.Range("X195").Value = 0.444
$ws.Range("Y195").Value = -1
$ws.Range("Z195").Value = -1
$ws.Range("AA195").Value = 0.925
$ws.Range("AB195").Value = -1
$ws.Range("AC195").Value = 1.025
$ws.Range("AD195").Value = -1

# Row 224
$ws.Range("B224").Value = 7920464
$ws.Range("E224").Value = 'Kifisias FC'
$ws.Range("F224").Value = 'Giannina'
$ws.Range("G224").Value = 2
$ws.Range("H224").Value = 3
$ws.Range("I224").Value = 1
$ws.Range("J224").Value = 1
$ws.Range("K224").Value = 'A'
$ws.Range("L224").Value = 1.571
$ws.Range("M224").Value = 3.8
$ws.Range("N224").Value = 6.5
$ws.Range("O224").Value = 1.4
$ws.Range("P224").Value = 4.75
$ws.Range("Q224").Value = 7.5
$ws.Range("R224").Value = -1.25
$ws.Range("S224").Value = 2.05
$ws.Range("T224").Value = 1.8
$ws.Range("U224").Value = 2.75
$ws.Range("V224").Value = 1.925
$ws.Range("W224").Value = 1.925
$ws.Range("X224").Value = -1
$ws.Range("Y224").Value = -1
$ws.Range("Z224").Value = 6.5
$ws.Range("AA224").Value = -1
$ws.Range("AB224").Value = 0.8
$ws.Range("AC224").Value = 0.925
$ws.Range("AD224").Value = -1

# Row 225
$ws.Range("B225").Value = 7920465
$ws.Range("E225").Value = 'Panetolikos'
$ws.Range("F225").Value = 'Atromitos Athinon'
$ws.Range("G225").Value = 1
$ws.Range("H225").Value = 0
$ws.Range("I225").Value = 0
$ws.Range("J225").Value = 0
$ws.Range("K225").Value = 'H'
$ws.Range("L225").Value = 1.7
$ws.Range("M225").Value = 3.6
$ws.Range("N225").Value = 5.25
$ws.Range("O225").Value = 1.666
$ws.Range("P225").Value = 3.75
$ws.Range("Q225").Value = 5.5
$ws.Range("R225").Value = -0.75
$ws.Range("S225").Value = 1.85
$ws.Range("T225").Value = 2
$ws.Range("U225").Value = 2.5
$ws.Range("V225").Value = 1.95
$ws.Range("W225").Value = 1.9
$ws.Range("X225").Value = 0.6659999999999999
$ws.Range("Y225").Value = -1
$ws.Range("Z225").Value = -1
$ws.Range("AA225").Value = 0.425
$ws.Range("AB225").Value = -0.5
$ws.Range("AC225").Value = -1
$ws.Range("AD225").Value = 0.8999999999999999

# Row 230
$ws.Range("B230").Value = 7920466
$ws.Range("E230").Value = 'Atromitos Athinon'
$ws.Range("F230").Value = 'Asteras Tripolis'
$ws.Range("G230").Value = 0
$ws.Range("H230").Value = 1
$ws.Range("I230").Value = 0
$ws.Range("J230").Value = 0
$ws.Range("K230").Value = 'A'
$ws.Range("L230").Value = 2.3
$ws.Range("M230").Value = 3.3
$ws.Range("N230").Value = 3.2
$ws.Range("O230").Value = 2.15
$ws.Range("P230").Value = 3.4
$ws.Range("Q230").Value = 3.3
$ws.Range("R230").Value = -0.25
$ws.Range("S230").Value = 1.85
$ws.Range("T230").Value = 2
$ws.Range("U230").Value = 2.5
$ws.Range("V230").Value = 1.975
$ws.Range("W230").Value = 1.875
$ws.Range("X230").Value = -1
$ws.Range("Y230").Value = -1
$ws.Range("Z230").Value = 2.3
$ws.Range("AA230").Value = -1
$ws.Range("AB230").Value = 1
$ws.Range("AC230").Value = -1
$ws.Range("AD230").Value = 0.875

# Row 232
$ws.Range("B232").Value = 7920467
$ws.Range("E232").Value = 'OFI Crete'
$ws.Range("F232").Value = 'Panetolikos'
$ws.Range("G232").Value = 1
$ws.Range("H232").Value = 2
$ws.Range("I232").Value = 0
$ws.Range("J232").Value = 1
$ws.Range("K232").Value = 'A'
$ws.Range("L232").Value = 2
$ws.Range("M232").Value = 3.5
$ws.Range("N232").Value = 3.75
$ws.Range("O232").Value = 2.05
$ws.Range("P232").Value = 3.4
$ws.Range("Q232").Value = 3.5
$ws.Range("R232").Value = -0.25
$ws.Range("S232").Value = 1.8
$ws.Range("T232").Value = 2.05
$ws.Range("U232").Value = 2.5
$ws.Range("V232").Value = 1.825
$ws.Range("W232").Value = 2.025
$ws.Range("X232").Value = -1
$ws.Range("Y232").Value = -1
$ws.Range("Z232").Value = 2.5
$ws.Range("AA232").Value = -1
$ws.Range("AB232").Value = 1.05
$ws.Range("AC232").Value = 0.825
$ws.Range("AD232").Value = -1
